$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formats from column M (2019) into new column N (2020), then fill values
$ws.Range("M4:M16").Copy()
$ws.Range("N4:N16").PasteSpecial(-4122)

$ws.Range("N4").Value = 2020
$ws.Range("N5").Value = 588.70000000000005
$ws.Range("N6").Value = 62.2
$ws.Range("N7").Value = 99.4
$ws.Range("N8").Value = 6.1
$ws.Range("N9").Value = "-"
$ws.Range("N10").Value = 71
$ws.Range("N11").Value = 136.30000000000001
$ws.Range("N12").Value = 103.3
$ws.Range("N13").Value = 103.2
$ws.Range("N14").Value = 1.8
$ws.Range("N15").Value = "-"
$ws.Range("N16").Value = 5.4

# N10 gets its own number format (0.0) distinct from the rest of the column
$ws.Range("N10").NumberFormat = "0.0"

$ws.Range("P15").Select() | Out-Null
